$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Programs" Cypher query (shown in column C for every tab row) gets two
# whitespace tweaks: the WHERE clause is indented two spaces, and a trailing
# space is added after the ['Melanoma'] literal.
$programsQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
  WHERE diag.disease_term IN ['Melanoma'] 
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$ws.Range("C2").Value = $programsQuery
$ws.Range("C3").Value = $programsQuery
$ws.Range("C4").Value = $programsQuery
$ws.Range("C5").Value = $programsQuery

# Keep the autofit row heights pinned to the values from the saved workbook.
$ws.Rows(2).RowHeight = 154.75
$ws.Rows(3).RowHeight = 156
$ws.Rows(4).RowHeight = 133.75
$ws.Rows(5).RowHeight = 174.65

# Selection moves from B5 to B2.
$ws.Range("B2").Select()
